$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0053680981595092
$ws.Range("C2").Value = 0.00230061349693252
$ws.Range("D2").Value = 0.0498466257668712
$ws.Range("E2").Value = 0.031441717791411
$ws.Range("F2").Value = 0.00230061349693252
$ws.Range("G2").Value = 0.00920245398773006
$ws.Range("H2").Value = 0.00920245398773006
$ws.Range("I2").Value = 0.00230061349693252
$ws.Range("J2").Value = 0.00920245398773006
$ws.Range("K2").Value = 0.996165644171779
$ws.Range("L2").Value = 0.000766871165644172
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0.000766871165644172
$ws.Range("O2").Value = 0.00306748466257669
$ws.Range("P2").Value = 0.00153374233128834
$ws.Range("Q2").Value = 0.00460122699386503
$ws.Range("R2").Value = 0.941717791411043
$ws.Range("S2").Value = 0.0107361963190184
$ws.Range("T2").Value = 0.00460122699386503
$ws.Range("U2").Value = 0.00996932515337423
$ws.Range("V2").Value = 0.849693251533742
$ws.Range("W2").Value = 0.159509202453988
$ws.Range("X2").Value = 0.0368098159509202
$ws.Range("B3").Value = 0.976993865030675
$ws.Range("C3").Value = 0.973159509202454
$ws.Range("D3").Value = 0.00613496932515337
$ws.Range("E3").Value = 0.953220858895706
$ws.Range("F3").Value = 0.0184049079754601
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0.980061349693252
$ws.Range("M3").Value = 0.00153374233128834
$ws.Range("N3").Value = 0.00460122699386503
$ws.Range("O3").Value = 0.00383435582822086
$ws.Range("P3").Value = 0.997699386503067
$ws.Range("Q3").Value = 0.977760736196319
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.000766871165644172
$ws.Range("T3").Value = 0.0191717791411043
$ws.Range("U3").Value = 0.0176380368098159
$ws.Range("V3").Value = 0.0145705521472393
$ws.Range("W3").Value = 0.000766871165644172
$ws.Range("X3").Value = 0.0237730061349693
$ws.Range("B4").Value = 0.0153374233128834
$ws.Range("C4").Value = 0.0191717791411043
$ws.Range("D4").Value = 0.943251533742331
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0.00230061349693252
$ws.Range("G4").Value = 0.989263803680982
$ws.Range("H4").Value = 0.989263803680982
$ws.Range("I4").Value = 0.996932515337423
$ws.Range("J4").Value = 0.990030674846626
$ws.Range("K4").Value = 0.00383435582822086
$ws.Range("L4").Value = 0.0153374233128834
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0145705521472393
$ws.Range("O4").Value = 0.0153374233128834
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0.0138036809815951
$ws.Range("R4").Value = 0.0575153374233129
$ws.Range("S4").Value = 0.983128834355828
$ws.Range("T4").Value = 0.975460122699387
$ws.Range("U4").Value = 0.968558282208589
$ws.Range("V4").Value = 0.131901840490798
$ws.Range("W4").Value = 0.825153374233129
$ws.Range("X4").Value = 0.93941717791411
$ws.Range("B5").Value = 0.00230061349693252
$ws.Range("C5").Value = 0.0053680981595092
$ws.Range("D5").Value = 0.000766871165644172
$ws.Range("E5").Value = 0.0153374233128834
$ws.Range("F5").Value = 0.976993865030675
$ws.Range("G5").Value = 0.000766871165644172
$ws.Range("H5").Value = 0.000766871165644172
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0.000766871165644172
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0.00383435582822086
$ws.Range("M5").Value = 0.998466257668712
$ws.Range("N5").Value = 0.980061349693252
$ws.Range("O5").Value = 0.977760736196319
$ws.Range("P5").Value = 0.000766871165644172
$ws.Range("Q5").Value = 0.00383435582822086
$ws.Range("R5").Value = 0.000766871165644172
$ws.Range("S5").Value = 0.0053680981595092
$ws.Range("T5").Value = 0.000766871165644172
$ws.Range("U5").Value = 0.00383435582822086
$ws.Range("V5").Value = 0.00383435582822086
$ws.Range("W5").Value = 0.0145705521472393
$ws.Range("X5").Value = 0
